$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fix the misspelled "Febuary" -> "February" date entries throughout the log.
# ---------------------------------------------------------------------------
$ws.Range("B275").Value = "February 3 2018"
$ws.Range("B276").Value = "February 4 2018"
$ws.Range("B278").Value = "February 5 2018"
$ws.Range("B280").Value = "February 6 2018"
$ws.Range("B283").Value = "February 7 2018"
$ws.Range("B285").Value = "February 9 2018"
$ws.Range("B288").Value = "February 10 2018"
$ws.Range("B292").Value = "February 11 2018"
$ws.Range("B295").Value = "February 12 2018"
$ws.Range("B297").Value = "February 13 2018"
$ws.Range("B301").Value = "February 14 2018"
$ws.Range("B302").Value = "February 15 2018"
$ws.Range("B303").Value = "February 16 2018"
$ws.Range("B305").Value = "February 17 2018"
$ws.Range("B308").Value = "February 18 2018"
$ws.Range("B311").Value = "February 19 2018"
$ws.Range("B313").Value = "February 20 2018"

# ---------------------------------------------------------------------------
# Populate the new rows in the Log table with the latest activity entries.
# ---------------------------------------------------------------------------
$ws.Range("B317").Value = "February 22 2018"
$ws.Range("C317").Value = 0.75
$ws.Range("D317").Value = 0.77083333333333337
$ws.Range("F317").Value = "Changed add journal button's color. Added showing hexagram detail modal feature to search reading page."

$ws.Range("B318").Value = "February 23 2018"
$ws.Range("C318").Value = 0.72916666666666663
$ws.Range("D318").Value = 0.75
$ws.Range("F318").Value = "Starting to extract some code relates to show hexagram detail modal."

$ws.Range("C319").Value = 0.86111111111111116
$ws.Range("D319").Value = 1.0326388888888889
$ws.Range("F319").Value = "Finished refactoring for HexagramDetailModal."

$ws.Range("B320").Value = "February 24 2018"
$ws.Range("C320").Value = 0.65763888888888888

# ---------------------------------------------------------------------------
# Update the view state (scroll position / active selection) to reflect the
# newly added rows at the bottom of the log.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 305
$win.ScrollColumn = 1
$ws.Range("C320").Select()
